$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2684
$ws.Range("F5").Value = 928
$ws.Range("F7").Value = 2200
$ws.Range("F8").Value = 1808
$ws.Range("F9").Value = 211
$ws.Range("F11").Value = 2464
$ws.Range("F12").Value = 540
$ws.Range("F13").Value = 228
$ws.Range("F17").Value = 109
$ws.Range("F18").Value = 9139
$ws.Range("F20").Value = 7084
$ws.Range("F21").Value = 11585
$ws.Range("F25").Value = 342
$ws.Range("F27").Value = 2556
$ws.Range("I27").Value = "//i0.hdslb.com/bfs/openplatform/202410/zOFCXyVt1728723765464.jpeg"
$ws.Range("F29").Value = 195
$ws.Range("F30").Value = 2492
$ws.Range("F31").Value = 674
$ws.Range("F33").Value = 4500
$ws.Range("F34").Value = 859
$ws.Range("F35").Value = 348
$ws.Range("F37").Value = 516

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 145

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2684
$ws.Range("F7").Value = 928
$ws.Range("F9").Value = 2200
$ws.Range("F11").Value = 1808
$ws.Range("F13").Value = 211
$ws.Range("F14").Value = 2464
$ws.Range("F16").Value = 540
$ws.Range("F17").Value = 228
$ws.Range("F21").Value = 109
$ws.Range("F22").Value = 9139
$ws.Range("F24").Value = 7084
$ws.Range("F25").Value = 11585
$ws.Range("F29").Value = 342
$ws.Range("F33").Value = 2556
$ws.Range("I33").Value = "//i0.hdslb.com/bfs/openplatform/202410/zOFCXyVt1728723765464.jpeg"
$ws.Range("F37").Value = 195
$ws.Range("F39").Value = 4500
$ws.Range("F46").Value = 516
